# Add a new "geno" column (K) derived from the genotype number embedded
# in column A (e.g. "IPa06-C-TF" -> "06", "OPa-01-C-TF" -> "01").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell K1: same text style as the other header cells (bold, centered).
$ws.Range("K1").Value = "geno"
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "geno"

# Keep the two-digit genotype codes ("01".."10") as text so leading zeros survive.
$ws.Range("K2:K294").NumberFormat = "@"

$genoValues = @(
  @(2, "06"),
  @(3, "07"),
  @(4, "08"),
  @(5, "09"),
  @(6, "10"),
  @(7, "06"),
  @(8, "07"),
  @(9, "08"),
  @(10, "09"),
  @(11, "10"),
  @(12, "06"),
  @(13, "07"),
  @(14, "08"),
  @(15, "09"),
  @(16, "10"),
  @(17, "06"),
  @(18, "07"),
  @(19, "08"),
  @(20, "09"),
  @(21, "10"),
  @(22, "06"),
  @(23, "07"),
  @(24, "08"),
  @(25, "09"),
  @(26, "10"),
  @(27, "06"),
  @(28, "07"),
  @(29, "08"),
  @(30, "09"),
  @(31, "10"),
  @(32, "06"),
  @(33, "07"),
  @(34, "08"),
  @(35, "09"),
  @(36, "10"),
  @(37, "06"),
  @(38, "07"),
  @(39, "08"),
  @(40, "09"),
  @(41, "10"),
  @(42, "06"),
  @(43, "07"),
  @(44, "08"),
  @(45, "09"),
  @(46, "10"),
  @(47, "06"),
  @(48, "07"),
  @(49, "08"),
  @(50, "09"),
  @(51, "10"),
  @(52, "06"),
  @(53, "07"),
  @(54, "08"),
  @(55, "09"),
  @(56, "10"),
  @(57, "06"),
  @(58, "07"),
  @(59, "08"),
  @(60, "09"),
  @(61, "10"),
  @(62, "06"),
  @(63, "07"),
  @(64, "08"),
  @(65, "09"),
  @(66, "10"),
  @(67, "06"),
  @(68, "07"),
  @(69, "08"),
  @(70, "09"),
  @(71, "10"),
  @(72, "06"),
  @(73, "07"),
  @(74, "08"),
  @(75, "09"),
  @(76, "10"),
  @(77, "06"),
  @(78, "07"),
  @(79, "08"),
  @(80, "09"),
  @(81, "10"),
  @(82, "06"),
  @(83, "07"),
  @(84, "08"),
  @(85, "09"),
  @(86, "10"),
  @(87, "06"),
  @(88, "07"),
  @(89, "08"),
  @(90, "09"),
  @(91, "10"),
  @(92, "06"),
  @(93, "07"),
  @(94, "08"),
  @(95, "09"),
  @(96, "10"),
  @(97, "06"),
  @(98, "07"),
  @(99, "08"),
  @(100, "09"),
  @(101, "10"),
  @(102, "06"),
  @(103, "07"),
  @(104, "08"),
  @(105, "09"),
  @(106, "10"),
  @(107, "06"),
  @(108, "07"),
  @(109, "08"),
  @(110, "09"),
  @(111, "10"),
  @(112, "01"),
  @(113, "02"),
  @(114, "03"),
  @(115, "04"),
  @(116, "05"),
  @(117, "01"),
  @(118, "02"),
  @(119, "03"),
  @(120, "04"),
  @(121, "05"),
  @(122, "01"),
  @(123, "02"),
  @(124, "03"),
  @(125, "04"),
  @(126, "05"),
  @(127, "01"),
  @(128, "02"),
  @(129, "03"),
  @(130, "04"),
  @(131, "05"),
  @(132, "01"),
  @(133, "02"),
  @(134, "03"),
  @(135, "04"),
  @(136, "05"),
  @(137, "01"),
  @(138, "02"),
  @(139, "03"),
  @(140, "04"),
  @(141, "05"),
  @(142, "01"),
  @(143, "02"),
  @(144, "03"),
  @(145, "04"),
  @(146, "05"),
  @(147, "01"),
  @(148, "02"),
  @(149, "03"),
  @(150, "04"),
  @(151, "05"),
  @(152, "01"),
  @(153, "02"),
  @(154, "03"),
  @(155, "04"),
  @(156, "05"),
  @(157, "01"),
  @(158, "02"),
  @(159, "03"),
  @(160, "04"),
  @(161, "05"),
  @(162, "01"),
  @(163, "02"),
  @(164, "03"),
  @(165, "04"),
  @(166, "05"),
  @(167, "01"),
  @(168, "02"),
  @(169, "03"),
  @(170, "04"),
  @(171, "05"),
  @(172, "01"),
  @(173, "02"),
  @(174, "03"),
  @(175, "04"),
  @(176, "05"),
  @(177, "01"),
  @(178, "02"),
  @(179, "03"),
  @(180, "04"),
  @(181, "05"),
  @(182, "01"),
  @(183, "02"),
  @(184, "03"),
  @(185, "04"),
  @(186, "05"),
  @(187, "01"),
  @(188, "02"),
  @(189, "03"),
  @(190, "04"),
  @(191, "05"),
  @(192, "01"),
  @(193, "02"),
  @(194, "03"),
  @(195, "04"),
  @(196, "05"),
  @(197, "01"),
  @(198, "02"),
  @(199, "03"),
  @(200, "04"),
  @(201, "05"),
  @(202, "01"),
  @(203, "02"),
  @(204, "03"),
  @(205, "04"),
  @(206, "05"),
  @(207, "01"),
  @(208, "02"),
  @(209, "03"),
  @(210, "04"),
  @(211, "05"),
  @(212, "01"),
  @(213, "02"),
  @(214, "03"),
  @(215, "04"),
  @(216, "05"),
  @(217, "01"),
  @(218, "02"),
  @(219, "03"),
  @(220, "04"),
  @(221, "05"),
  @(222, "06"),
  @(223, "05"),
  @(224, "09"),
  @(225, "04"),
  @(226, "08"),
  @(227, "10"),
  @(228, "10"),
  @(229, "09"),
  @(230, "09"),
  @(231, "10"),
  @(232, "08"),
  @(233, "06"),
  @(234, "09"),
  @(235, "07"),
  @(236, "10"),
  @(237, "08"),
  @(238, "10"),
  @(239, "09"),
  @(240, "09"),
  @(241, "07"),
  @(242, "06"),
  @(243, "10"),
  @(244, "09"),
  @(245, "10"),
  @(246, "07"),
  @(247, "09"),
  @(248, "09"),
  @(249, "07"),
  @(250, "07"),
  @(251, "08"),
  @(252, "10"),
  @(253, "08"),
  @(254, "06"),
  @(255, "07"),
  @(256, "07"),
  @(257, "07"),
  @(258, "08"),
  @(259, "08"),
  @(260, "06"),
  @(261, "07"),
  @(262, "08"),
  @(263, "10"),
  @(264, "09"),
  @(265, "09"),
  @(266, "09"),
  @(267, "02"),
  @(268, "05"),
  @(269, "02"),
  @(270, "03"),
  @(271, "03"),
  @(272, "01"),
  @(273, "04"),
  @(274, "01"),
  @(275, "04"),
  @(276, "05"),
  @(277, "05"),
  @(278, "01"),
  @(279, "02"),
  @(280, "03"),
  @(281, "04"),
  @(282, "05"),
  @(283, "02"),
  @(284, "03"),
  @(285, "02"),
  @(286, "02"),
  @(287, "03"),
  @(288, "01"),
  @(289, "02"),
  @(290, "03"),
  @(291, "04"),
  @(292, "05"),
  @(293, "01"),
  @(294, "01")
)

foreach ($pair in $genoValues) {
    $rowNum = $pair[0]
    $genoVal = $pair[1]
    $ws.Cells.Item($rowNum, 11).Value = $genoVal
}
